$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the Price column as Text so numeric-looking strings
# (e.g. "0.9993", "24.94") are stored as literal text, matching the
# original inlineStr cell type, instead of being parsed as numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value2 = '29.376.88'
$ws.Range("E2").Value2 = '  +0.02%  '

$ws.Range("D3").Value2 = '1.841.25'
$ws.Range("E3").Value2 = '  -0.28%  '

$ws.Range("D4").Value2 = '0.9993'
$ws.Range("E4").Value2 = '  +0.12%  '

$ws.Range("D5").Value2 = '239.18'
$ws.Range("E5").Value2 = '  -0.40%  '

$ws.Range("D6").Value2 = '0.6265'
$ws.Range("E6").Value2 = '  -0.19%  '

$ws.Range("E7").Value2 = '  +0.10%  '

$ws.Range("D8").Value2 = '0.07433'
$ws.Range("E8").Value2 = '  -0.81%  '

$ws.Range("D9").Value2 = '0.2892'
$ws.Range("E9").Value2 = '  -0.31%  '

$ws.Range("D10").Value2 = '24.94'
$ws.Range("E10").Value2 = '  +1.83%  '

$ws.Range("D11").Value2 = '0.07720'
$ws.Range("E11").Value2 = '  -0.26%  '

$ws.Range("D12").Value2 = '1.828.91'
$ws.Range("E12").Value2 = '  -0.95%  '

$ws.Range("D13").Value2 = '4.976'
$ws.Range("E13").Value2 = '  -0.36%  '

$ws.Range("D14").Value2 = '0.6737'
$ws.Range("E14").Value2 = '  -0.99%  '

$ws.Range("D15").Value2 = '0.00001028'
$ws.Range("E15").Value2 = '  -1.90%  '

$ws.Range("D16").Value2 = '81.77'
$ws.Range("E16").Value2 = '  -0.33%  '

$ws.Range("D17").Value2 = '6.206'
$ws.Range("E17").Value2 = '  +0.11%  '

$ws.Range("D18").Value2 = '29.423.31'
$ws.Range("E18").Value2 = '  +0.13%  '

$ws.Range("D19").Value2 = '234.20'
$ws.Range("E19").Value2 = '  +2.28%  '

$ws.Range("D20").Value2 = '12.32'
$ws.Range("E20").Value2 = '  -0.15%  '

$ws.Range("E21").Value2 = '  +0.17%  '

$ws.Range("D22").Value2 = '7.293'
$ws.Range("E22").Value2 = '  -2.63%  '

$ws.Range("E23").Value2 = '  +0.13%  '

$ws.Range("D24").Value2 = '157.93'
$ws.Range("E24").Value2 = '  -0.40%  '

$ws.Range("D25").Value2 = '8.500'
$ws.Range("E25").Value2 = '  +0.77%  '

$ws.Range("D26").Value2 = '0.1343'
$ws.Range("E26").Value2 = '  -1.91%  '

$ws.Range("E27").Value2 = '  -1.35%  '

$ws.Range("D28").Value2 = '0.07250'
$ws.Range("E28").Value2 = '  +11.54%  '

$ws.Range("D29").Value2 = '1.467'
$ws.Range("E29").Value2 = '  +4.08%  '

$ws.Range("E30").Value2 = '  -0.13%  '

$ws.Range("D31").Value2 = '4.036'
$ws.Range("E31").Value2 = '  -1.05%  '

$ws.Range("E32").Value2 = '  -1.63%  '

$ws.Range("D33").Value2 = '1.816'
$ws.Range("E33").Value2 = '  -0.80%  '

$ws.Range("D34").Value2 = '1.138'
$ws.Range("E34").Value2 = '  -0.35%  '

$ws.Range("D35").Value2 = '0.6965'
$ws.Range("E35").Value2 = '  -0.17%  '

$ws.Range("E36").Value2 = '  -0.20%  '

$ws.Range("D37").Value2 = '0.01843'
$ws.Range("E37").Value2 = '  +0.43%  '

$ws.Range("D38").Value2 = '6.930'
$ws.Range("E38").Value2 = '  +2.65%  '

$ws.Range("D39").Value2 = '2.817'
$ws.Range("E39").Value2 = '  -0.52%  '

$ws.Range("D40").Value2 = '1.235.41'
$ws.Range("E40").Value2 = '  -2.40%  '

$ws.Range("D41").Value2 = '0.9604'
$ws.Range("E41").Value2 = '  +4.65%  '

$ws.Range("D42").Value2 = '1.001'
$ws.Range("E42").Value2 = '  +0.16%  '

$ws.Range("D43").Value2 = '1.999.18'
$ws.Range("E43").Value2 = '  -0.45%  '

$ws.Range("D44").Value2 = '100.83'
$ws.Range("E44").Value2 = '  -0.61%  '

$ws.Range("D45").Value2 = '65.37'
$ws.Range("E45").Value2 = '  -1.33%  '

$ws.Range("E46").Value2 = '  +0.94%  '

$ws.Range("D47").Value2 = '1.718'
$ws.Range("E47").Value2 = '  -0.28%  '

$ws.Range("D48").Value2 = '6.949'
$ws.Range("E48").Value2 = '  -1.94%  '

$ws.Range("D49").Value2 = '8.882'
$ws.Range("E49").Value2 = '  -0.94%  '

$ws.Range("D50").Value2 = '0.1131'
$ws.Range("E50").Value2 = '  -2.77%  '

$ws.Range("D51").Value2 = '0.3898'
$ws.Range("E51").Value2 = '  -1.56%  '

# Restore the column to the default (unstyled) look now that the
# values are committed as text, so no stray style index lingers.
$priceRange.Style = "Normal"
